$d = $word.ActiveDocument

# The "CRONOGRAMA DE HITOS DEL PROYECTO" table is the 5th table in the document.
$t = $d.Tables.Item(5)

# Row 5 = "Hito 3: ..." row. Shrink its height (418 -> 255 twentieths of a point).
$row3 = $t.Rows.Item(5)
$row3.Height = 12.75

# Update the Hito 3 title.
$cell3title = $row3.Cells.Item(1)
$cell3title.Range.Find.Execute("Hito 3:  Desarrollo Front-End y Back-End", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Hito 3:  Desarrollo y Pruebas", 2)

# Update the Hito 3 date (scoped to this cell so it can't clobber the "24/10/2025" in the next row).
$cell3date = $row3.Cells.Item(2)
$cell3date.Range.Find.Execute("4/10/2025", $true, $false, $false, $false, $false, `
    $true, 1, $false, "24/10/2025", 2)

# The old "Hito 4: Implementación de base de datos" and "Hito 5: Pruebas de calidad (QA)" rows
# are removed entirely.
$t.Rows.Item(6).Delete()
$t.Rows.Item(6).Delete()

# The remaining last row (formerly "Hito 6: Entrega final y capacitación") becomes the new
# "Hito 4" row; only its title text changes, the date (13/11/2025) stays as-is.
$lastRow = $t.Rows.Item($t.Rows.Count)
$cellLastTitle = $lastRow.Cells.Item(1)
$cellLastTitle.Range.Find.Execute("Hito 6: Entrega final y capacitación", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Hito 4: Despliegue", 2)
